$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 302000000.0
$ws.Range("C4").Value = 319000000.0
$ws.Range("D4").Value = 339000000.0
$ws.Range("E4").Value = 327000000.0
$ws.Range("F4").Value = 312000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 143000000.0
$ws.Range("C14").Value = 131000000.0
$ws.Range("D14").Value = 129000000.0
$ws.Range("E14").Value = 145000000.0
$ws.Range("F14").Value = 159000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -210000000.0
$ws.Range("C22").Value = -147000000.0
$ws.Range("D22").Value = -145000000.0
$ws.Range("E22").Value = -156000000.0
$ws.Range("F22").Value = -157000000.0
